# "change state to simple Column"
#
# 1. Duplicate the "Status" sheet to the end of the workbook, renaming the
#    copy "Status_1" - this preserves the original 3-column (migration /
#    stateId / stateName) layout on the new copy.
# 2. On the original "Status" sheet, delete the now-redundant "stateId"
#    column (column B), collapsing it down to a simple 2-column
#    (migration / stateName) layout.
# 3. Stamp a default page setup on "Status" (paper size / orientation),
#    matching what Excel records after the sheet is touched.
# 4. On "Permissions", add a stray "0,3" value down at F24, copying the
#    header formatting.
# 5. Make "Permissions" the active/selected sheet (instead of "API"), and
#    restore each sheet's last-used cell selection.

$wb = $excel.ActiveWorkbook

# --- 1. Copy "Status" -> "Status_1" (placed after the last sheet) -------
$wsStatus = $wb.Worksheets.Item("Status")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsStatus.Copy($null, $lastSheet)

$wsStatus1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsStatus1.Name = "Status_1"
$wsStatus1.Activate()
$wsStatus1.Range("B22").Select()

# --- 2. Reduce "Status" to a simple 2-column layout ----------------------
$wsStatus.Columns("B").Delete()

$wsStatus.PageSetup.PaperSize = 9
$wsStatus.PageSetup.Orientation = 1

$wsStatus.Activate()
$wsStatus.Range("N11").Select()

# --- 4. Add stray "0,3" value on "Permissions" ---------------------------
$wsPerm = $wb.Worksheets.Item("Permissions")
$wsPerm.Range("A1").Copy()
$wsPerm.Range("F24").PasteSpecial(-4122)
$wsPerm.Range("F24").Value = "0,3"

# --- 5. Make "Permissions" the active sheet, restore selections ---------
$wsPerm.Activate()
$wsPerm.Range("F24").Select()

$wsApi = $wb.Worksheets.Item("API")
$wsApi.Activate()
$wsApi.Range("B16").Select()

$wsPerm.Activate()
